$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new "record" columns (AD, AE, AF).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting used by the rest of the header row (bold, border,
# centered) by copying the formats from the neighboring header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins/Losses/Ties) for every data row (2-42).
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 75
    $ws.Cells.Item($r, 32).Value = 0
}
